$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (A1:C1) ---
$ws.Range("A1").Value = "nombre"
$ws.Range("B1").Value = "edad"
$ws.Range("C1").Value = "ciudad"

# New column C needs the same width as A/B (stored width=15) and the
# same header styling as A1/B1 (fill/border/bold/center).
$ws.Columns(3).ColumnWidth = 14.17
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# --- Update first data row (row 2) ---
$ws.Range("A2").Value = "Erick"
$ws.Range("B2").Value = 28
$ws.Range("C2").Value = "Mérida"

# --- Remove the other data rows (old rows 3-5: Marketing/Operaciones/Finanzas) ---
$ws.Rows("3:5").Delete()

# --- Point the chart series at the single remaining data row ---
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection(1)
$ser.XValues = "='Reporte'!`$A`$2"
$ser.Values = "='Reporte'!`$B`$2"
